# Natmi following Dr Hou advice
# Replace the LR-pair result rows with the updated NATMI computation (15 sending/target cluster combinations).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 15,20
$data[0,0] = "ECs"
$data[0,1] = "Cxcl1"
$data[0,2] = "Xcr1"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 19.79062366666667
$data[0,7] = 59.371871
$data[0,8] = 0.04690933339254189
$data[0,9] = 0.04982137012915518
$data[0,10] = 1
$data[0,11] = 0.3333333333333333
$data[0,12] = 0.01506066666666667
$data[0,13] = 0.045182
$data[0,14] = 0.008628232818297613
$data[0,15] = 0.008628232818297613
$data[0,16] = 0.2980599861691111
$data[0,17] = 2.682539875522
$data[0,18] = 0.0004047446498619941
$data[0,19] = 0.0004298703808009291

$data[1,0] = "ECs"
$data[1,1] = "Cxcl1"
$data[1,2] = "Xcr1"
$data[1,3] = "M1"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 19.79062366666667
$data[1,7] = 59.371871
$data[1,8] = 0.04690933339254189
$data[1,9] = 0.04982137012915518
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 1.034921333333333
$data[1,13] = 3.104764
$data[1,14] = 0.592904843474591
$data[1,15] = 0.5929048434745909
$data[1,16] = 20.48173863260489
$data[1,17] = 184.335647693444
$data[1,18] = 0.02781277097260246
$data[1,19] = 0.02953933165811641

$data[2,0] = "ECs"
$data[2,1] = "Cxcl1"
$data[2,2] = "Xcr1"
$data[2,3] = "M2"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 19.79062366666667
$data[2,7] = 59.371871
$data[2,8] = 0.04690933339254189
$data[2,9] = 0.04982137012915518
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 0.695528
$data[2,13] = 2.086584
$data[2,14] = 0.3984669237071115
$data[2,15] = 0.3984669237071114
$data[2,16] = 13.76493289762933
$data[2,17] = 123.884396078664
$data[2,18] = 0.01869181777007745
$data[2,19] = 0.01985216809023784

$data[3,0] = "FAPs"
$data[3,1] = "Cxcl1"
$data[3,2] = "Xcr1"
$data[3,3] = "ECs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 260.6166636666667
$data[3,7] = 781.849991
$data[3,8] = 0.617734648968278
$data[3,9] = 0.6560823691590862
$data[3,10] = 1
$data[3,11] = 0.3333333333333333
$data[3,12] = 0.01506066666666667
$data[3,13] = 0.045182
$data[3,14] = 0.008628232818297613
$data[3,15] = 0.008628232818297613
$data[3,16] = 3.925060699262445
$data[3,17] = 35.325546293362
$data[3,18] = 0.005329958371227652
$data[3,19] = 0.005660831429084878

$data[4,0] = "FAPs"
$data[4,1] = "Cxcl1"
$data[4,2] = "Xcr1"
$data[4,3] = "M1"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 260.6166636666667
$data[4,7] = 781.849991
$data[4,8] = 0.617734648968278
$data[4,9] = 0.6560823691590862
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 1.034921333333333
$data[4,13] = 3.104764
$data[4,14] = 0.592904843474591
$data[4,15] = 0.5929048434745909
$data[4,16] = 269.7177450507916
$data[4,17] = 2427.459705457124
$data[4,18] = 0.3662578653553683
$data[4,19] = 0.3889944143927068

$data[5,0] = "FAPs"
$data[5,1] = "Cxcl1"
$data[5,2] = "Xcr1"
$data[5,3] = "M2"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 260.6166636666667
$data[5,7] = 781.849991
$data[5,8] = 0.617734648968278
$data[5,9] = 0.6560823691590862
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 0.695528
$data[5,13] = 2.086584
$data[5,14] = 0.3984669237071115
$data[5,15] = 0.3984669237071114
$data[5,16] = 181.2661868467493
$data[5,17] = 1631.395681620744
$data[5,18] = 0.2461468252416821
$data[5,19] = 0.2614271233372945

$data[6,0] = "M1"
$data[6,1] = "Cxcl1"
$data[6,2] = "Xcr1"
$data[6,3] = "ECs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 19.92354333333333
$data[6,7] = 59.77063
$data[6,8] = 0.04722439031359255
$data[6,9] = 0.05015598514796319
$data[6,10] = 1
$data[6,11] = 0.3333333333333333
$data[6,12] = 0.01506066666666667
$data[6,13] = 0.045182
$data[6,14] = 0.008628232818297613
$data[6,15] = 0.008628232818297613
$data[6,16] = 0.3000618449622222
$data[6,17] = 2.70055660466
$data[6,18] = 0.0004074630343278351
$data[6,19] = 0.0004327575170877037

$data[7,0] = "M1"
$data[7,1] = "Cxcl1"
$data[7,2] = "Xcr1"
$data[7,3] = "M1"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 19.92354333333333
$data[7,7] = 59.77063
$data[7,8] = 0.04722439031359255
$data[7,9] = 0.05015598514796319
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 1.034921333333333
$data[7,13] = 3.104764
$data[7,14] = 0.592904843474591
$data[7,15] = 0.5929048434745909
$data[7,16] = 20.61930003125778
$data[7,17] = 185.57370028132
$data[7,18] = 0.02799956974706358
$data[7,19] = 0.02973772652346702

$data[8,0] = "M1"
$data[8,1] = "Cxcl1"
$data[8,2] = "Xcr1"
$data[8,3] = "M2"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 19.92354333333333
$data[8,7] = 59.77063
$data[8,8] = 0.04722439031359255
$data[8,9] = 0.05015598514796319
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 0.695528
$data[8,13] = 2.086584
$data[8,14] = 0.3984669237071115
$data[8,15] = 0.3984669237071114
$data[8,16] = 13.85738224754667
$data[8,17] = 124.71644022792
$data[8,18] = 0.01881735753220113
$data[8,19] = 0.01998550110740846

$data[9,0] = "M2"
$data[9,1] = "Cxcl1"
$data[9,2] = "Xcr1"
$data[9,3] = "ECs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 47.58211666666667
$data[9,7] = 142.74635
$data[9,8] = 0.1127829729792156
$data[9,9] = 0.1197843123039854
$data[9,10] = 1
$data[9,11] = 0.3333333333333333
$data[9,12] = 0.01506066666666667
$data[9,13] = 0.045182
$data[9,14] = 0.008628232818297613
$data[9,15] = 0.008628232818297613
$data[9,16] = 0.7166183984111112
$data[9,17] = 6.4495655857
$data[9,18] = 0.0009731177488044409
$data[9,19] = 0.001033526934538457

$data[10,0] = "M2"
$data[10,1] = "Cxcl1"
$data[10,2] = "Xcr1"
$data[10,3] = "M1"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 47.58211666666667
$data[10,7] = 142.74635
$data[10,8] = 0.1127829729792156
$data[10,9] = 0.1197843123039854
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 1.034921333333333
$data[10,13] = 3.104764
$data[10,14] = 0.592904843474591
$data[10,15] = 0.5929048434745909
$data[10,16] = 49.24374762348889
$data[10,17] = 443.1937286114
$data[10,18] = 0.06686957094084085
$data[10,19] = 0.07102069893730595

$data[11,0] = "M2"
$data[11,1] = "Cxcl1"
$data[11,2] = "Xcr1"
$data[11,3] = "M2"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 47.58211666666667
$data[11,7] = 142.74635
$data[11,8] = 0.1127829729792156
$data[11,9] = 0.1197843123039854
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 0.695528
$data[11,13] = 2.086584
$data[11,14] = 0.3984669237071115
$data[11,15] = 0.3984669237071114
$data[11,16] = 33.09469444093334
$data[11,17] = 297.8522499684
$data[11,18] = 0.04494028428957031
$data[11,19] = 0.04773008643214094

$data[12,0] = "sCs"
$data[12,1] = "Cxcl1"
$data[12,2] = "Xcr1"
$data[12,3] = "ECs"
$data[12,4] = 2
$data[12,5] = 1
$data[12,6] = 73.978012
$data[12,7] = 147.956024
$data[12,8] = 0.1753486543463721
$data[12,9] = 0.1241559632598098
$data[12,10] = 1
$data[12,11] = 0.3333333333333333
$data[12,12] = 0.01506066666666667
$data[12,13] = 0.045182
$data[12,14] = 0.008628232818297613
$data[12,15] = 0.008628232818297613
$data[12,16] = 1.114158179394667
$data[12,17] = 6.684949076368
$data[12,18] = 0.001512949014075693
$data[12,19] = 0.001071246556785644

$data[13,0] = "sCs"
$data[13,1] = "Cxcl1"
$data[13,2] = "Xcr1"
$data[13,3] = "M1"
$data[13,4] = 2
$data[13,5] = 1
$data[13,6] = 73.978012
$data[13,7] = 147.956024
$data[13,8] = 0.1753486543463721
$data[13,9] = 0.1241559632598098
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 1.034921333333333
$data[13,13] = 3.104764
$data[13,14] = 0.592904843474591
$data[13,15] = 0.5929048434745909
$data[13,16] = 76.56142281638934
$data[13,17] = 459.368536898336
$data[13,18] = 0.1039650664587159
$data[13,19] = 0.0736126719629946

$data[14,0] = "sCs"
$data[14,1] = "Cxcl1"
$data[14,2] = "Xcr1"
$data[14,3] = "M2"
$data[14,4] = 2
$data[14,5] = 1
$data[14,6] = 73.978012
$data[14,7] = 147.956024
$data[14,8] = 0.1753486543463721
$data[14,9] = 0.1241559632598098
$data[14,10] = 3
$data[14,11] = 1
$data[14,12] = 0.695528
$data[14,13] = 2.086584
$data[14,14] = 0.3984669237071115
$data[14,15] = 0.3984669237071114
$data[14,16] = 51.453778730336
$data[14,17] = 308.7226723820161
$data[14,18] = 0.06987063887358053
$data[14,19] = 0.04947204474002957

$ws.Range("A2:T16").Value = $data
